# Ave_SOOT_Mass_FRAC.xlsx edit:
#  - Rename Sheet1 headers to generic Var1/Var2/Var3
#  - Update the "Uncertainty" (col C) values on Sheet1 with recomputed figures
#    (including a previously-missing C8 value)
#  - Add a new "Mass_Frac_Soot" worksheet at the end with descriptive headers
#    (Position_cm / ave_Y_Soot / u_ave_Y_Soot) and the same data set

$wb = $excel.ActiveWorkbook

# ---- Sheet1: header rename + corrected uncertainty column ----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").Value = "Var1"
$ws1.Range("B1").Value = "Var2"
$ws1.Range("C1").Value = "Var3"

$ws1.Range("C2").Value = [double]"1.4821872093528022E-5"
$ws1.Range("C3").Value = [double]"2.8962459490505354E-4"
$ws1.Range("C4").Value = [double]"7.9435861828233061E-5"
$ws1.Range("C5").Value = [double]"6.8587613385216614E-5"
$ws1.Range("C6").Value = [double]"5.3561238416670722E-5"
$ws1.Range("C7").Value = [double]"4.0016422727251008E-5"
$ws1.Range("C8").Value = [double]"4.1300003580276173E-5"
$ws1.Range("C9").Value = [double]"5.5592372236095103E-5"
$ws1.Range("C10").Value = [double]"4.3615910695519825E-5"

$ws1.Columns.Item(1).ColumnWidth = 4.1666666
$ws1.Columns.Item(2).ColumnWidth = 11.1666666
$ws1.Columns.Item(3).ColumnWidth = 11.1666666

# ---- New worksheet: Mass_Frac_Soot (appended after Sheet3) ----
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Mass_Frac_Soot"

$newSheet.Range("A1").Value = "Position_cm"
$newSheet.Range("B1").Value = "ave_Y_Soot"
$newSheet.Range("C1").Value = "u_ave_Y_Soot"

$positions = @(2, 4, 6, 10, 14, 20, 30, 45, 60, 100)
$aveY = @([double]"3.6528022450646957E-4", [double]"6.3937065295494521E-4", [double]"1.0403569827979284E-3", [double]"1.3347755453614568E-3", [double]"1.3957587341662063E-3", [double]"1.1402221192398716E-3", [double]"9.4558349928821435E-4", [double]"3.745553362714471E-4", [double]"1.1743408933365021E-4", [double]"1.2500387622552602E-5")
$uY    = @([double]"1.4821872093528022E-5", [double]"2.8962459490505354E-4", [double]"7.9435861828233061E-5", [double]"6.8587613385216614E-5", [double]"5.3561238416670722E-5", [double]"4.0016422727251008E-5", [double]"4.1300003580276173E-5", [double]"5.5592372236095103E-5", [double]"4.3615910695519825E-5", [double]"1.3312392218280024E-5")

for ($i = 0; $i -lt $positions.Length; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 1).Value = $positions[$i]
    $newSheet.Cells.Item($row, 2).Value = $aveY[$i]
    $newSheet.Cells.Item($row, 3).Value = $uY[$i]
}

$newSheet.Columns.Item(1).ColumnWidth = 11.0
$newSheet.Columns.Item(2).ColumnWidth = 11.1666666
$newSheet.Columns.Item(3).ColumnWidth = 12.5766666

# Sheet1 stays the on-screen/selected tab in the saved file.
$ws1.Activate()

Write-Host "Edit complete"
